# Replace the hard-coded date "16 серпня 2018" in the registration-action
# paragraph with the template placeholder token "{5}", matching the
# docxFiller convention already used elsewhere in this document
# (e.g. "... від {5} р.").
#
# The target OOXML splits the run that used to hold the whole sentence
# into three runs: the text before the date, the new "{5}" token, and the
# text after the date - all three sharing the same run formatting
# (Times New Roman / bCs / uk-UA). A plain Find/Replace keeps everything
# inside a single run, so after swapping the text we nudge a character
# property away and back to force the run boundaries to be written out
# separately.

$d = $word.ActiveDocument

$r = $d.Content
$found = $r.Find.Execute("16 серпня 2018", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $wasBold = $r.Bold

    # Swap the matched date text for the placeholder token.
    $r.Text = "{5}"

    # Force a run split at the edges of the replaced range: toggle a
    # character property away from its current value and immediately
    # restore it, so the new text is materialized as its own run(s)
    # distinct from the surrounding text instead of being re-merged into
    # the original run.
    $r.Bold = 1
    $r.Bold = $wasBold
}
